$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Localization")

# --- Add new localization rows for Start / Quit / Apply / Cancel ---

# Insert 3 blank rows for Start / Quit / Cancel
$ws.Rows("8:10").Insert()

# Row 8: Start
$ws.Range("A8").Value = "Start"
$ws.Range("B8").Value = "Start"
$ws.Range("C8").Value = "Start"

# Row 9: Quit (German filled in later)
$ws.Range("A9").Value = "Quit"
$ws.Range("B9").Value = "Quit"

# Row 10: Cancel, temporarily here (German filled in later)
$ws.Range("A10").Value = "Cancel"
$ws.Range("B10").Value = "Cancel"

# Insert the Apply row above Cancel, pushing Cancel down to row 11
$ws.Rows("10:10").Insert()
$ws.Range("A10").Value = "Apply"
$ws.Range("B10").Value = "Apply"

# Fill in the German translations
$ws.Range("C10").Value = "Anwenden"
$ws.Range("C9").Value = "Beenden"
$ws.Range("C11").Value = "Abbrechen"

# Column D is always the empty-string formula on data rows
$ws.Range("D8").Formula = "="""""
$ws.Range("D9").Formula = "="""""
$ws.Range("D10").Formula = "="""""
$ws.Range("D11").Formula = "="""""

# --- Add the new "Select Player" row at the bottom ---
$ws.Range("B15").Value = "Select Player"
$ws.Range("A15").Value = "SelectPlayer"
$ws.Range("C15").Value = "Spieler Ändern"
$ws.Range("D15").Formula = "="""""

# Move the active selection down to the row below the new data
$ws.Range("A16:XFD16").Select()
